$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.674.97"
$ws.Range("E2").Value2 = "  +4.16%  "
$ws.Range("D3").Value2 = "3.499.45"
$ws.Range("E3").Value2 = "  +2.21%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "591.82"
$ws.Range("E5").Value2 = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "169.19"
$ws.Range("E6").Value2 = "  +3.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.999"
$ws.Range("E7").Value2 = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.600"
$ws.Range("E8").Value2 = "  +8.12%  "
$ws.Range("D9").Value2 = "3.496.62"
$ws.Range("E9").Value2 = "  +2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.129"
$ws.Range("E10").Value2 = "  +7.28%  "
$ws.Range("E11").Value2 = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.441"
$ws.Range("E12").Value2 = "  +3.97%  "
$ws.Range("D13").Value2 = "4.103.65"
$ws.Range("E13").Value2 = "  +2.28%  "
$ws.Range("E14").Value2 = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "28.24"
$ws.Range("E15").Value2 = "  +4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.0000180"
$ws.Range("E16").Value2 = "  +3.48%  "
$ws.Range("D17").Value2 = "66.646.01"
$ws.Range("E17").Value2 = "  +4.03%  "
$ws.Range("D18").Value2 = "3.514.24"
$ws.Range("E18").Value2 = "  +3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.34"
$ws.Range("E19").Value2 = "  +3.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "14.19"
$ws.Range("E20").Value2 = "  +3.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "393.11"
$ws.Range("E21").Value2 = "  +3.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "7.98"
$ws.Range("E22").Value2 = "  +2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "73.30"
$ws.Range("E23").Value2 = "  +3.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.999"
$ws.Range("E24").Value2 = "  +0.00%  "
$ws.Range("E25").Value2 = "  +3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.0000123"
$ws.Range("E26").Value2 = "  +5.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.27"
$ws.Range("E27").Value2 = "  +7.43%  "
$ws.Range("E28").Value2 = "  +2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.997"
$ws.Range("E29").Value2 = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "6.39"
$ws.Range("E30").Value2 = "  +4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.48"
$ws.Range("E31").Value2 = "  +5.30%  "
$ws.Range("E32").Value2 = "  +3.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "23.64"
$ws.Range("E33").Value2 = "  +2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "7.45"
$ws.Range("E34").Value2 = "  +4.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.00"
$ws.Range("E35").Value2 = "  +0.03%  "
$ws.Range("E36").Value2 = "  +8.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "162.69"
$ws.Range("E37").Value2 = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.883"
$ws.Range("E38").Value2 = "  +2.73%  "
$ws.Range("E39").Value2 = "  +5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "27.68"
$ws.Range("E40").Value2 = "  +4.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "6.80"
$ws.Range("E41").Value2 = "  +4.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "4.68"
$ws.Range("E42").Value2 = "  +5.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0748"
$ws.Range("E43").Value2 = "  +2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "26.52"
$ws.Range("E44").Value2 = "  +2.14%  "
$ws.Range("D45").Value2 = "2.792.00"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "43.21"
$ws.Range("E46").Value2 = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0311"
$ws.Range("E47").Value2 = "  +1.35%  "
$ws.Range("E48").Value2 = "  +2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "351.38"
$ws.Range("E49").Value2 = "  +6.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.11"
$ws.Range("E50").Value2 = "  +5.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "33.76"
$ws.Range("E51").Value2 = "  +12.23%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
